$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 2.597803197262351
$ws.Cells.Item(2, 3).Value = 0.7119485641435404
$ws.Cells.Item(2, 4).Value = 0.03291285572721492
$ws.Cells.Item(2, 5).Value = 1.361407321067915
$ws.Cells.Item(2, 6).Value = 0.6798212621317248
$ws.Cells.Item(2, 7).Value = 0.0007871586232472039
$ws.Cells.Item(2, 8).Value = 0.01199277220940109
$ws.Cells.Item(2, 9).Value = 0.005767783202876409
$ws.Cells.Item(2, 16).Value = 0.7161795053628026
$ws.Cells.Item(2, 17).Value = 2.014744168414211
$ws.Cells.Item(3, 2).Value = 2.26066291172782
$ws.Cells.Item(3, 3).Value = 0.6281112901472454
$ws.Cells.Item(3, 4).Value = 0.02978638338481687
$ws.Cells.Item(3, 5).Value = 1.18640414660436
$ws.Cells.Item(3, 6).Value = 0.6162561042387509
$ws.Cells.Item(3, 7).Value = 0.0007907025278933496
$ws.Cells.Item(3, 8).Value = 0.008635475661823488
$ws.Cells.Item(3, 9).Value = 0.003802278604460696
$ws.Cells.Item(3, 16).Value = 0.7252262043596858
$ws.Cells.Item(3, 17).Value = 1.84802363948063
$ws.Cells.Item(4, 2).Value = 2.053033621404154
$ws.Cells.Item(4, 3).Value = 0.5769465878343851
$ws.Cells.Item(4, 4).Value = 0.02786294938359291
$ws.Cells.Item(4, 5).Value = 1.079243740448391
$ws.Cells.Item(4, 6).Value = 0.5777840511071304
$ws.Cells.Item(4, 7).Value = 0.0007929509702298748
$ws.Cells.Item(4, 8).Value = 0.006788786149038328
$ws.Cells.Item(4, 9).Value = 0.002808140480463628
$ws.Cells.Item(4, 16).Value = 0.7314904468061201
$ws.Cells.Item(4, 17).Value = 1.747272484144673
$ws.Cells.Item(5, 2).Value = 1.966867226263588
$ws.Cells.Item(5, 3).Value = 0.5570060540242423
$ws.Cells.Item(5, 4).Value = 0.02711732349066409
$ws.Cells.Item(5, 5).Value = 1.03560817179681
$ws.Cells.Item(5, 6).Value = 0.561636489855097
$ws.Cells.Item(5, 7).Value = 0.0007938913891312019
$ws.Cells.Item(5, 8).Value = 0.006084113324819063
$ws.Cells.Item(5, 9).Value = 0.002507873075721534
$ws.Cells.Item(5, 16).Value = 0.734578658242711
$ws.Cells.Item(5, 17).Value = 1.704511633347835
$ws.Cells.Item(6, 2).Value = 1.950852273096388
$ws.Cells.Item(6, 3).Value = 0.5547293813975784
$ws.Cells.Item(6, 4).Value = 0.02704163866897247
$ws.Cells.Item(6, 5).Value = 1.028341575679136
$ws.Cells.Item(6, 6).Value = 0.5582333237960952
$ws.Cells.Item(6, 7).Value = 0.0007940555458227334
$ws.Cells.Item(6, 8).Value = 0.005967871774951172
$ws.Cells.Item(6, 9).Value = 0.00253600607557658
$ws.Cells.Item(6, 16).Value = 0.7355487681603776
$ws.Cells.Item(6, 17).Value = 1.694901505474832
$ws.Cells.Item(7, 2).Value = 2.047220876175572
$ws.Cells.Item(7, 3).Value = 0.579499765070608
$ws.Cells.Item(7, 4).Value = 0.02798496047145704
$ws.Cells.Item(7, 5).Value = 1.078590135761104
$ws.Cells.Item(7, 6).Value = 0.5755708414071776
$ws.Cells.Item(7, 7).Value = 0.0007929820356354663
$ws.Cells.Item(7, 8).Value = 0.006772998661915941
$ws.Cells.Item(7, 9).Value = 0.002993613051261335
$ws.Cells.Item(7, 16).Value = 0.7327682632237398
$ws.Cells.Item(7, 17).Value = 1.739771223154349
$ws.Cells.Item(8, 2).Value = 2.475449915545028
$ws.Cells.Item(8, 3).Value = 0.6867448461070467
$ws.Cells.Item(8, 4).Value = 0.03201250899027031
$ws.Cells.Item(8, 5).Value = 1.300903148091805
$ws.Cells.Item(8, 6).Value = 0.655126001160383
$ws.Cells.Item(8, 7).Value = 0.0007883887738010078
$ws.Cells.Item(8, 8).Value = 0.01077826887721478
$ws.Cells.Item(8, 9).Value = 0.00526474178485703
$ws.Cells.Item(8, 16).Value = 0.7208571998837741
$ws.Cells.Item(8, 17).Value = 1.947695453102057
$ws.Cells.Item(9, 2).Value = 3.31937571960276
$ws.Cells.Item(9, 3).Value = 0.89486404983694
$ws.Cells.Item(9, 4).Value = 0.03964780683899249
$ws.Cells.Item(9, 5).Value = 1.74031629218878
$ws.Cells.Item(9, 6).Value = 0.819821406831295
$ws.Cells.Item(9, 7).Value = 0.0007799115522071909
$ws.Cells.Item(9, 8).Value = 0.02055750911132859
$ws.Cells.Item(9, 9).Value = 0.01132324856504763
$ws.Cells.Item(9, 16).Value = 0.7010876780096567
$ws.Cells.Item(9, 17).Value = 2.382696616817753
$ws.Cells.Item(10, 2).Value = 3.917712953324951
$ws.Cells.Item(10, 3).Value = 1.045199349042434
$ws.Cells.Item(10, 4).Value = 0.04618919566138402
$ws.Cells.Item(10, 5).Value = 1.963957886527425
$ws.Cells.Item(10, 6).Value = 0.93232782730243
$ws.Cells.Item(10, 7).Value = 0.0007741951350906243
$ws.Cells.Item(10, 8).Value = 0.02843873462482227
$ws.Cells.Item(10, 9).Value = 0.01710822360147279
$ws.Cells.Item(10, 16).Value = 0.696941688860889
$ws.Cells.Item(10, 17).Value = 2.671041276445521
$ws.Cells.Item(11, 2).Value = 4.026737869601959
$ws.Cells.Item(11, 3).Value = 1.062850731216088
$ws.Cells.Item(11, 4).Value = 0.05755877076850879
$ws.Cells.Item(11, 5).Value = 1.273281716279527
$ws.Cells.Item(11, 6).Value = 0.8817242950077855
$ws.Cells.Item(11, 7).Value = 0.0007731507858902821
$ws.Cells.Item(11, 8).Value = 0.04364912973931112
$ws.Cells.Item(11, 9).Value = 0.01881957031504289
$ws.Cells.Item(11, 16).Value = 0.7518756670614692
$ws.Cells.Item(11, 17).Value = 2.455349917182559
$ws.Cells.Item(12, 2).Value = 3.998808706553689
$ws.Cells.Item(12, 3).Value = 1.04164719867569
$ws.Cells.Item(12, 4).Value = 0.06634448787927028
$ws.Cells.Item(12, 5).Value = 0.7757246460136571
$ws.Cells.Item(12, 6).Value = 0.8193973257164942
$ws.Cells.Item(12, 7).Value = 0.0007732720225684575
$ws.Cells.Item(12, 8).Value = 0.07933126607416341
$ws.Cells.Item(12, 9).Value = 0.01876606120865443
$ws.Cells.Item(12, 16).Value = 0.8018015336147926
$ws.Cells.Item(12, 17).Value = 2.230018738977265
$ws.Cells.Item(13, 2).Value = 3.856906087453297
$ws.Cells.Item(13, 3).Value = 0.9941221507705222
$ws.Cells.Item(13, 4).Value = 0.07374033173330474
$ws.Cells.Item(13, 5).Value = 0.3981850151505242
$ws.Cells.Item(13, 6).Value = 0.7415224755416716
$ws.Cells.Item(13, 7).Value = 0.0007743164300377484
$ws.Cells.Item(13, 8).Value = 0.1320877733521399
$ws.Cells.Item(13, 9).Value = 0.01762405812788259
$ws.Cells.Item(13, 16).Value = 0.8516298066913066
$ws.Cells.Item(13, 17).Value = 1.974162583525896
$ws.Cells.Item(14, 2).Value = 3.705445098923747
$ws.Cells.Item(14, 3).Value = 0.9502030679651909
$ws.Cells.Item(14, 4).Value = 0.07837039433519521
$ws.Cells.Item(14, 5).Value = 0.2067866589870491
$ws.Cells.Item(14, 6).Value = 0.6801878403898911
$ws.Cells.Item(14, 7).Value = 0.0007754502749329242
$ws.Cells.Item(14, 8).Value = 0.1794594729324785
$ws.Cells.Item(14, 9).Value = 0.01645342251417592
$ws.Cells.Item(14, 16).Value = 0.8871022934417567
$ws.Cells.Item(14, 17).Value = 1.782298678331358
$ws.Cells.Item(15, 2).Value = 3.642840219142954
$ws.Cells.Item(15, 3).Value = 0.9350401674155364
$ws.Cells.Item(15, 4).Value = 0.07915881600541752
$ws.Cells.Item(15, 5).Value = 0.1691515420980636
$ws.Cells.Item(15, 6).Value = 0.6612674466339854
$ws.Cells.Item(15, 7).Value = 0.0007759511119717162
$ws.Cells.Item(15, 8).Value = 0.1913068405711869
$ws.Cells.Item(15, 9).Value = 0.01602642478418481
$ws.Cells.Item(15, 16).Value = 0.8957326648395139
$ws.Cells.Item(15, 17).Value = 1.726423683432074
$ws.Cells.Item(16, 2).Value = 3.415509897697746
$ws.Cells.Item(16, 3).Value = 0.8831420113135096
$ws.Cells.Item(16, 4).Value = 0.0746227652443352
$ws.Cells.Item(16, 5).Value = 0.1649757372158618
$ws.Cells.Item(16, 6).Value = 0.6290727610731253
$ws.Cells.Item(16, 7).Value = 0.000778154564283884
$ws.Cells.Item(16, 8).Value = 0.1761469891596477
$ws.Cells.Item(16, 9).Value = 0.01389576357840028
$ws.Cells.Item(16, 16).Value = 0.8846610237793584
$ws.Cells.Item(16, 17).Value = 1.655776269964292
$ws.Cells.Item(17, 2).Value = 3.321024665012828
$ws.Cells.Item(17, 3).Value = 0.8665758524583111
$ws.Cells.Item(17, 4).Value = 0.068686150943293
$ws.Cells.Item(17, 5).Value = 0.2495924249046269
$ws.Cells.Item(17, 6).Value = 0.6365869619460085
$ws.Cells.Item(17, 7).Value = 0.0007792633737742622
$ws.Cells.Item(17, 8).Value = 0.1374730480102357
$ws.Cells.Item(17, 9).Value = 0.01290158642736738
$ws.Cells.Item(17, 16).Value = 0.8577170160333907
$ws.Cells.Item(17, 17).Value = 1.70188707370113
$ws.Cells.Item(18, 2).Value = 3.333655881952268
$ws.Cells.Item(18, 3).Value = 0.8764144046545539
$ws.Cells.Item(18, 4).Value = 0.06103533650605186
$ws.Cells.Item(18, 5).Value = 0.4836990774361354
$ws.Cells.Item(18, 6).Value = 0.6812400819957105
$ws.Cells.Item(18, 7).Value = 0.000779451795486122
$ws.Cells.Item(18, 8).Value = 0.08540656275748404
$ws.Cells.Item(18, 9).Value = 0.0125143417278295
$ws.Cells.Item(18, 16).Value = 0.8143454827812491
$ws.Cells.Item(18, 17).Value = 1.86245199742703
$ws.Cells.Item(19, 2).Value = 3.423425903977261
$ws.Cells.Item(19, 3).Value = 0.9129069093284556
$ws.Cells.Item(19, 4).Value = 0.05334773736532838
$ws.Cells.Item(19, 5).Value = 0.9192880719533036
$ws.Cells.Item(19, 6).Value = 0.7509527051255276
$ws.Cells.Item(19, 7).Value = 0.0007787894674839167
$ws.Cells.Item(19, 8).Value = 0.04300193645821793
$ws.Cells.Item(19, 9).Value = 0.01316157143171637
$ws.Cells.Item(19, 16).Value = 0.7672507512976097
$ws.Cells.Item(19, 17).Value = 2.100504158686675
$ws.Cells.Item(20, 2).Value = 3.744961384763428
$ws.Cells.Item(20, 3).Value = 1.014395815254744
$ws.Cells.Item(20, 4).Value = 0.04493882850729136
$ws.Cells.Item(20, 5).Value = 1.900556197826361
$ws.Cells.Item(20, 6).Value = 0.8957133485068596
$ws.Cells.Item(20, 7).Value = 0.0007757310825716802
$ws.Cells.Item(20, 8).Value = 0.02618642432198026
$ws.Cells.Item(20, 9).Value = 0.01606431092290084
$ws.Cells.Item(20, 16).Value = 0.7024029010301973
$ws.Cells.Item(20, 17).Value = 2.571054206805428
$ws.Cells.Item(21, 2).Value = 4.234385973037206
$ws.Cells.Item(21, 3).Value = 1.13884384872128
$ws.Cells.Item(21, 4).Value = 0.04853039978854667
$ws.Cells.Item(21, 5).Value = 2.232452845818756
$ws.Cells.Item(21, 6).Value = 1.002870546799215
$ws.Cells.Item(21, 7).Value = 0.0007711437311507585
$ws.Cells.Item(21, 8).Value = 0.03398707486947439
$ws.Cells.Item(21, 9).Value = 0.02118666147670556
$ws.Cells.Item(21, 16).Value = 0.691096490122753
$ws.Cells.Item(21, 17).Value = 2.862923816794051
$ws.Cells.Item(22, 2).Value = 4.551661220263384
$ws.Cells.Item(22, 3).Value = 1.213745153735829
$ws.Cells.Item(22, 4).Value = 0.05114237283794409
$ws.Cells.Item(22, 5).Value = 2.398415709812426
$ws.Cells.Item(22, 6).Value = 1.070752524721385
$ws.Cells.Item(22, 7).Value = 0.000768271936691179
$ws.Cells.Item(22, 8).Value = 0.03907437156944749
$ws.Cells.Item(22, 9).Value = 0.02456729336757402
$ws.Cells.Item(22, 16).Value = 0.6855652115969377
$ws.Cells.Item(22, 17).Value = 3.045853247365073
$ws.Cells.Item(23, 2).Value = 4.388043389984887
$ws.Cells.Item(23, 3).Value = 1.170230267484214
$ws.Cells.Item(23, 4).Value = 0.04958402825251795
$ws.Cells.Item(23, 5).Value = 2.309816508383619
$ws.Cells.Item(23, 6).Value = 1.03681460987049
$ws.Cells.Item(23, 7).Value = 0.0007697827481974752
$ws.Cells.Item(23, 8).Value = 0.03633694744138571
$ws.Cells.Item(23, 9).Value = 0.02250664103468569
$ws.Cells.Item(23, 16).Value = 0.6867202211229895
$ws.Cells.Item(23, 17).Value = 2.95620379021554
$ws.Cells.Item(24, 2).Value = 3.759689213248521
$ws.Cells.Item(24, 3).Value = 1.012181686639536
$ws.Cells.Item(24, 4).Value = 0.04394261903344443
$ws.Cells.Item(24, 5).Value = 1.976461631646146
$ws.Cells.Item(24, 6).Value = 0.9067763043870798
$ws.Cells.Item(24, 7).Value = 0.0007756338485891714
$ws.Cells.Item(24, 8).Value = 0.02669420700939273
$ws.Cells.Item(24, 9).Value = 0.01574906392351672
$ws.Cells.Item(24, 16).Value = 0.6955573263841259
$ws.Cells.Item(24, 17).Value = 2.610298883614547
$ws.Cells.Item(25, 2).Value = 3.082796153344248
$ws.Cells.Item(25, 3).Value = 0.8434117984034515
$ws.Cells.Item(25, 4).Value = 0.03782965901152835
$ws.Cells.Item(25, 5).Value = 1.620805740526038
$ws.Cells.Item(25, 6).Value = 0.7709755821200588
$ws.Cells.Item(25, 7).Value = 0.000782176244044853
$ws.Cells.Item(25, 8).Value = 0.01766107065060085
$ws.Cells.Item(25, 9).Value = 0.009773803154761396
$ws.Cells.Item(25, 16).Value = 0.7082226143536374
$ws.Cells.Item(25, 17).Value = 2.250509085723536
